$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.511.54"
$ws.Range("E2").Value = "  +0.22%  "
$ws.Range("D3").Value = "1.869.36"
$ws.Range("E3").Value = "  -0.72%  "
$ws.Range("E4").Value = "  -0.56%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.006"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.85%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5051"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.92%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3890"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.27%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08346"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.25%  "
$ws.Range("E10").Value = "  -0.53%  "
$ws.Range("E11").Value = "  -2.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.191"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.63%  "
$ws.Range("D13").Value = "1.870.92"
$ws.Range("E13").Value = "  -0.85%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.40"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.23%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.221"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.68%  "
$ws.Range("E16").Value = "  -0.74%  "
$ws.Range("E17").Value = "  -0.91%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "90.92"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.58%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06689"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.96%  "
$ws.Range("E20").Value = "  -0.90%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.006"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.73%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.901"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.93%  "
$ws.Range("D23").Value = "28.551.46"
$ws.Range("E23").Value = "  +0.26%  "
$ws.Range("E24").Value = "  -1.41%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.227"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.67%  "
$ws.Range("D26").Value = "2.085.41"
$ws.Range("E26").Value = "  -0.77%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "161.50"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.15%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.59"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.17%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.335"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.91%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.55"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.44%  "
$ws.Range("E31").Value = "  -2.67%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.037"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.33%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.765"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.54%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.604"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.14%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02443"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.13%  "
$ws.Range("E36").Value = "  -0.52%  "
$ws.Range("E37").Value = "  -1.89%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.825"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.84%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.038"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.91%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.250"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.91%  "
$ws.Range("E41").Value = "  -0.77%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6405"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.44%  "
$ws.Range("E43").Value = "  -1.32%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.004"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6013"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.17%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.92"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.74%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.689"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.65%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.997"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.65%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.211"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.83%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "121.61"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.57%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.176"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -8.61%  "
